$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column C ("Bisnis") : fill top to bottom first ---
$ws.Range("C3").Value = "Distribusi"
$ws.Range("C4").Value = "Pergudangan"
$ws.Range("C5").Value = "All"
$ws.Range("C6").Value = "Distribusi"
$ws.Range("C7").Value = "Transportasi"
$ws.Range("C8").Value = "Transportasi"
$ws.Range("C9").Value = "All"
$ws.Range("C11").Value = "Transportasi"
$ws.Range("C12").Value = "All"
$ws.Range("C13").Value = "Distribusi"
$ws.Range("C14").Value = "Distribusi"
$ws.Range("C15").Value = "Transportasi"
$ws.Range("C16").Value = "All"
$ws.Range("C17").Value = "All"
$ws.Range("C18").Value = "All"
$ws.Range("C19").Value = "All"
$ws.Range("C20").Value = "All"
$ws.Range("C21").Value = "All"
$ws.Range("C22").Value = "All"
$ws.Range("C23").Value = "All"
$ws.Range("C24").Value = "All"
$ws.Range("C25").Value = "All"
$ws.Range("C26").Value = "All"
$ws.Range("C27").Value = "Pergudangan"
$ws.Range("C28").Value = "All"
$ws.Range("C29").Value = "All"
$ws.Range("C30").Value = "Transportasi"
$ws.Range("C31").Value = "Transportasi"
$ws.Range("C32").Value = "Transportasi"
$ws.Range("C33").Value = "Transportasi"
$ws.Range("C34").Value = "Transportasi"
$ws.Range("C35").Value = "Transportasi"
$ws.Range("C36").Value = "Transportasi"
$ws.Range("C37").Value = "Transportasi"
$ws.Range("C38").Value = "Distribusi"
$ws.Range("C39").Value = "Distribusi"
$ws.Range("C40").Value = "Distribusi"
$ws.Range("C41").Value = "Distribusi"
$ws.Range("C42").Value = "Pergudangan"
$ws.Range("C43").Value = "Pergudangan"
$ws.Range("C44").Value = "Pergudangan"
$ws.Range("C45").Value = "Pergudangan"

# --- Column H ("Aspect") : fill top to bottom afterwards ---
$ws.Range("H3").Value = "Produktifitas"
$ws.Range("H4").Value = "Produktifitas"
$ws.Range("H5").Value = "Kualitas"
$ws.Range("H6").Value = "Kualitas"
$ws.Range("H7").Value = "Produktifitas"
$ws.Range("H8").Value = "Kualitas"
$ws.Range("H9").Value = "Produktifitas"
$ws.Range("H10").Value = "Produktifitas"
$ws.Range("H11").Value = "Produktifitas"
$ws.Range("H12").Value = "Keuangan/Produktivitas "
$ws.Range("H13").Value = "Kualitas"
$ws.Range("H14").Value = "Kualitas"
$ws.Range("H15").Value = "Produktifitas"
$ws.Range("H16").Value = "Produktifitas"
$ws.Range("H17").Value = "Kualitas"
$ws.Range("H18").Value = "Kualitas"
$ws.Range("H19").Value = "Keuangan"
$ws.Range("H20").Value = "Keuangan"
$ws.Range("H21").Value = "Keuangan"
$ws.Range("H22").Value = "Kualitas"
$ws.Range("H23").Value = "Produktifitas"
$ws.Range("H24").Value = "Keuangan"
$ws.Range("H25").Value = "Keuangan"
$ws.Range("H26").Value = "Kualitas"
$ws.Range("H27").Value = "Produktifitas"
$ws.Range("H28").Value = "Keuangan"
$ws.Range("H29").Value = "Keuangan"
$ws.Range("H30").Value = "Kualitas"
$ws.Range("H31").Value = "Produktifitas"
$ws.Range("H32").Value = "Produktifitas"
$ws.Range("H33").Value = "Produktifitas"
$ws.Range("H34").Value = "Keuangan/Produktivitas "
$ws.Range("H35").Value = "Produktifitas"
$ws.Range("H36").Value = "Kualitas"
$ws.Range("H37").Value = "Kualitas"
$ws.Range("H38").Value = "Kualitas"
$ws.Range("H39").Value = "Kualitas"
$ws.Range("H40").Value = "Kualitas"
$ws.Range("H41").Value = "Produktifitas"
$ws.Range("H42").Value = "Keuangan/Produktivitas "
$ws.Range("H43").Value = "Produktifitas"
$ws.Range("H44").Value = "Produktifitas"
$ws.Range("H45").Value = "Produktifitas"

# Row height adjustments observed in the target file
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(12).RowHeight = 28.8

# View changes: scroll so column G is the leftmost visible column, and select H5
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("H5").Select()
